# Auto-generated edit applying updated crypto Price (D) / Volume(1h) (E) values
# as scraped on Tue Oct  3 17:31:45 UTC 2023 with GitHub Actions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.430.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.654.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.02%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("E6").Value = '  -1.46%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.19'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("E9").Value = '  -1.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0615'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0879'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.888.27'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.649.48'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.28%  '
$ws.Range("E14").Value = '  -2.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.571'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.88'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.436.40'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.98%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '233.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0727'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.51%  '
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("E22").Value = '  -3.01%  '
$ws.Range("E23").Value = '  -2.42%  '
$ws.Range("E24").Value = '  -1.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.79%  '
$ws.Range("E26").Value = '  -1.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("E29").Value = '  -2.09%  '
$ws.Range("E30").Value = '  -1.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.20'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.28%  '
$ws.Range("E32").Value = '  -2.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.464.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.10'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.50%  '
$ws.Range("E35").Value = '  -4.05%  '
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("E37").Value = '  -3.54%  '
$ws.Range("E38").Value = '  -3.31%  '
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("E42").Value = '  -0.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.74'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.31%  '
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.797.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.783'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.59%  '
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0106'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.00%  '
$ws.Range("E50").Value = '  -1.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.76'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.14%  '
